$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.777.72"
$ws.Range("E2").Value = "  +0.48%  "

$ws.Range("D3").Value = "2.529.52"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.63"
$ws.Range("E5").Value = "  +0.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.97"
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("E7").Value = "  -1.76%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -1.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.11"
$ws.Range("E10").Value = "  -0.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0811"
$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.56"
$ws.Range("E12").Value = "  -2.51%  "

$ws.Range("E13").Value = "  -3.32%  "

$ws.Range("D14").Value = "2.917.35"
$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("D15").Value = "2.528.27"
$ws.Range("E15").Value = "  +0.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.21"
$ws.Range("E16").Value = "  -2.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.853"
$ws.Range("E17").Value = "  -1.08%  "

$ws.Range("D18").Value = "42.864.51"
$ws.Range("E18").Value = "  +0.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.82"
$ws.Range("E19").Value = "  +4.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.86"
$ws.Range("E20").Value = "  -0.15%  "

$ws.Range("E21").Value = "  -0.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.95"
$ws.Range("E22").Value = "  -2.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.32"
$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.94"
$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.08"
$ws.Range("E25").Value = "  +1.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.79"
$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("E27").Value = "  -0.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.43"
$ws.Range("E28").Value = "  +3.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.56"
$ws.Range("E29").Value = "  +7.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.34"

$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.93"
$ws.Range("E32").Value = "  +1.73%  "

$ws.Range("E33").Value = "  +4.66%  "

$ws.Range("E34").Value = "  +0.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.11"
$ws.Range("E35").Value = "  -1.94%  "

$ws.Range("E36").Value = "  +1.91%  "

$ws.Range("E37").Value = "  -0.54%  "

$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("E39").Value = "  -1.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.42"
$ws.Range("E40").Value = "  -2.62%  "

$ws.Range("E41").Value = "  +14.60%  "

$ws.Range("E42").Value = "  -0.74%  "

$ws.Range("E43").Value = "  +0.38%  "

$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("E45").Value = "  -2.20%  "

$ws.Range("D46").Value = "2.036.44"
$ws.Range("E46").Value = "  +0.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.98"
$ws.Range("E47").Value = "  +0.84%  "

$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "107.15"
$ws.Range("E49").Value = "  +5.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.04"
$ws.Range("E50").Value = "  +2.52%  "

$ws.Range("D51").Value = "2.772.65"
$ws.Range("E51").Value = "  +0.55%  "
